$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.8203074518761176
$ws.Range("J2").Value = 0.8725723693674974
$ws.Range("M2").Value = 9.24193
$ws.Range("N2").Value = 27.72579
$ws.Range("O2").Value = 0.1468938537243544
$ws.Range("P2").Value = 0.1569651396557324
$ws.Range("Q2").Value = 0.5045077167699999
$ws.Range("R2").Value = 4.54056945093
$ws.Range("S2").Value = 0.1204981228448883
$ws.Range("T2").Value = 0.1369634438175025
$ws.Range("I3").Value = 0.8203074518761176
$ws.Range("J3").Value = 0.8725723693674974
$ws.Range("O3").Value = 0.469548954544906
$ws.Range("P3").Value = 0.5017420086455576
$ws.Range("S3").Value = 0.3851745064338268
$ws.Range("T3").Value = 0.4378062132950616
$ws.Range("I4").Value = 0.8203074518761176
$ws.Range("J4").Value = 0.8725723693674974
$ws.Range("M4").Value = 7.349831333333333
$ws.Range("N4").Value = 22.049494
$ws.Range("O4").Value = 0.1168203014713749
$ws.Range("P4").Value = 0.1248296948454213
$ws.Range("Q4").Value = 0.4012199426553333
$ws.Range("R4").Value = 3.610979483898
$ws.Range("S4").Value = 0.09582856382738339
$ws.Range("T4").Value = 0.1089229425986909
$ws.Range("I5").Value = 0.8203074518761176
$ws.Range("J5").Value = 0.8725723693674974
$ws.Range("M5").Value = 12.1104985
$ws.Range("N5").Value = 24.220997
$ws.Range("O5").Value = 0.1924876941491673
$ws.Range("P5").Value = 0.1371233128688515
$ws.Range("Q5").Value = 0.6611000026164999
$ws.Range("R5").Value = 3.966600015699
$ws.Range("S5").Value = 0.1578990899050129
$ws.Range("T5").Value = 0.1196500140054944
$ws.Range("I6").Value = 0.8203074518761176
$ws.Range("J6").Value = 0.8725723693674974
$ws.Range("M6").Value = 4.671440333333334
$ws.Range("N6").Value = 14.014321
$ws.Range("O6").Value = 0.07424919611019735
$ws.Range("P6").Value = 0.079339843984437
$ws.Range("Q6").Value = 0.2550092563563334
$ws.Range("R6").Value = 2.295083307207
$ws.Range("S6").Value = 0.06090716886500614
$ws.Range("T6").Value = 0.06922975565074778
$ws.Range("G7").Value = 0.011958
$ws.Range("H7").Value = 0.023916
$ws.Range("I7").Value = 0.1796925481238824
$ws.Range("J7").Value = 0.1274276306325027
$ws.Range("M7").Value = 9.24193
$ws.Range("N7").Value = 27.72579
$ws.Range("O7").Value = 0.1468938537243544
$ws.Range("P7").Value = 0.1569651396557324
$ws.Range("Q7").Value = 0.11051499894
$ws.Range("R7").Value = 0.66308999364
$ws.Range("S7").Value = 0.0263957308794661
$ws.Range("T7").Value = 0.02000169583822986
$ws.Range("G8").Value = 0.011958
$ws.Range("H8").Value = 0.023916
$ws.Range("I8").Value = 0.1796925481238824
$ws.Range("J8").Value = 0.1274276306325027
$ws.Range("O8").Value = 0.469548954544906
$ws.Range("P8").Value = 0.5017420086455576
$ws.Range("Q8").Value = 0.353263263902
$ws.Range("R8").Value = 2.119579583412
$ws.Range("S8").Value = 0.08437444811107919
$ws.Range("T8").Value = 0.06393579535049608
$ws.Range("G9").Value = 0.011958
$ws.Range("H9").Value = 0.023916
$ws.Range("I9").Value = 0.1796925481238824
$ws.Range("J9").Value = 0.1274276306325027
$ws.Range("M9").Value = 7.349831333333333
$ws.Range("N9").Value = 22.049494
$ws.Range("O9").Value = 0.1168203014713749
$ws.Range("P9").Value = 0.1248296948454213
$ws.Range("Q9").Value = 0.087889283084
$ws.Range("R9").Value = 0.527335698504
$ws.Range("S9").Value = 0.02099173764399147
$ws.Range("T9").Value = 0.01590675224673037
$ws.Range("G10").Value = 0.011958
$ws.Range("H10").Value = 0.023916
$ws.Range("I10").Value = 0.1796925481238824
$ws.Range("J10").Value = 0.1274276306325027
$ws.Range("M10").Value = 12.1104985
$ws.Range("N10").Value = 24.220997
$ws.Range("O10").Value = 0.1924876941491673
$ws.Range("P10").Value = 0.1371233128688515
$ws.Range("Q10").Value = 0.144817341063
$ws.Range("R10").Value = 0.5792693642519999
$ws.Range("S10").Value = 0.03458860424415439
$ws.Range("T10").Value = 0.01747329886335711
$ws.Range("G11").Value = 0.011958
$ws.Range("H11").Value = 0.023916
$ws.Range("I11").Value = 0.1796925481238824
$ws.Range("J11").Value = 0.1274276306325027
$ws.Range("M11").Value = 4.671440333333334
$ws.Range("N11").Value = 14.014321
$ws.Range("O11").Value = 0.07424919611019735
$ws.Range("P11").Value = 0.079339843984437
$ws.Range("Q11").Value = 0.055861083506
$ws.Range("R11").Value = 0.335166501036
$ws.Range("S11").Value = 0.01334202724519122
$ws.Range("T11").Value = 0.01011008833368923
